$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$theme = $nm.Theme
$tcs = $theme.ThemeColorScheme
$c3 = $tcs.Item(3)
$c3.RGB = 123456
